$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column width for new column CD (82) to match existing data columns (width=12).
# ColumnWidth has a fixed +0.8333333333333334 offset baked into this runtime's
# char-width -> stored-width conversion, so we back it out to land exactly on 12.
$ws.Range("CD1").ColumnWidth = 11.166666666666666

# Row 1 header: new date label "2024/11/29" as literal text (not an auto-converted date).
$ws.Range("CD1").NumberFormat = "@"
$ws.Range("CD1").Value = "2024/11/29"
$ws.Range("B1").Copy()
$ws.Range("CD1").PasteSpecial(-4122)

# Data rows 2-53: new numeric column CD, each cell carries its own highlight style
# (1 = default, 2 = yellow fill, 3 = light-blue fill) inherited via a same-workbook
# format-only paste from a reference cell of that style (row 2 has all three: A2/D2/N2).
$ws.Range("CD2").Value = 208.8
$ws.Range("A2").Copy()
$ws.Range("CD2").PasteSpecial(-4122)
$ws.Range("CD3").Value = 154.6
$ws.Range("A2").Copy()
$ws.Range("CD3").PasteSpecial(-4122)
$ws.Range("CD4").Value = 161.3
$ws.Range("A2").Copy()
$ws.Range("CD4").PasteSpecial(-4122)
$ws.Range("CD5").Value = 149.2
$ws.Range("A2").Copy()
$ws.Range("CD5").PasteSpecial(-4122)
$ws.Range("CD6").Value = 134.5
$ws.Range("N2").Copy()
$ws.Range("CD6").PasteSpecial(-4122)
$ws.Range("CD7").Value = 306.9
$ws.Range("A2").Copy()
$ws.Range("CD7").PasteSpecial(-4122)
$ws.Range("CD8").Value = 147.5
$ws.Range("A2").Copy()
$ws.Range("CD8").PasteSpecial(-4122)
$ws.Range("CD9").Value = 148.3
$ws.Range("A2").Copy()
$ws.Range("CD9").PasteSpecial(-4122)
$ws.Range("CD10").Value = 131.4
$ws.Range("N2").Copy()
$ws.Range("CD10").PasteSpecial(-4122)
$ws.Range("CD11").Value = 168
$ws.Range("A2").Copy()
$ws.Range("CD11").PasteSpecial(-4122)
$ws.Range("CD12").Value = 155
$ws.Range("A2").Copy()
$ws.Range("CD12").PasteSpecial(-4122)
$ws.Range("CD13").Value = 134.2
$ws.Range("N2").Copy()
$ws.Range("CD13").PasteSpecial(-4122)
$ws.Range("CD14").Value = 160.7
$ws.Range("A2").Copy()
$ws.Range("CD14").PasteSpecial(-4122)
$ws.Range("CD15").Value = 186.2
$ws.Range("A2").Copy()
$ws.Range("CD15").PasteSpecial(-4122)
$ws.Range("CD16").Value = 144.2
$ws.Range("A2").Copy()
$ws.Range("CD16").PasteSpecial(-4122)
$ws.Range("CD17").Value = 124.6
$ws.Range("D2").Copy()
$ws.Range("CD17").PasteSpecial(-4122)
$ws.Range("CD18").Value = 173.2
$ws.Range("A2").Copy()
$ws.Range("CD18").PasteSpecial(-4122)
$ws.Range("CD19").Value = 144.3
$ws.Range("A2").Copy()
$ws.Range("CD19").PasteSpecial(-4122)
$ws.Range("CD20").Value = 145.6
$ws.Range("A2").Copy()
$ws.Range("CD20").PasteSpecial(-4122)
$ws.Range("CD21").Value = 133
$ws.Range("N2").Copy()
$ws.Range("CD21").PasteSpecial(-4122)
$ws.Range("CD22").Value = 166.6
$ws.Range("A2").Copy()
$ws.Range("CD22").PasteSpecial(-4122)
$ws.Range("CD23").Value = 165.6
$ws.Range("A2").Copy()
$ws.Range("CD23").PasteSpecial(-4122)
$ws.Range("CD24").Value = 140.6
$ws.Range("A2").Copy()
$ws.Range("CD24").PasteSpecial(-4122)
$ws.Range("CD25").Value = 148.5
$ws.Range("A2").Copy()
$ws.Range("CD25").PasteSpecial(-4122)
$ws.Range("CD26").Value = 178.1
$ws.Range("A2").Copy()
$ws.Range("CD26").PasteSpecial(-4122)
$ws.Range("CD27").Value = 159.3
$ws.Range("A2").Copy()
$ws.Range("CD27").PasteSpecial(-4122)
$ws.Range("CD28").Value = 134
$ws.Range("N2").Copy()
$ws.Range("CD28").PasteSpecial(-4122)
$ws.Range("CD29").Value = 133.5
$ws.Range("N2").Copy()
$ws.Range("CD29").PasteSpecial(-4122)
$ws.Range("CD30").Value = 113.5
$ws.Range("D2").Copy()
$ws.Range("CD30").PasteSpecial(-4122)
$ws.Range("CD31").Value = 124.6
$ws.Range("D2").Copy()
$ws.Range("CD31").PasteSpecial(-4122)
$ws.Range("CD32").Value = 166.9
$ws.Range("A2").Copy()
$ws.Range("CD32").PasteSpecial(-4122)
$ws.Range("CD33").Value = 211.5
$ws.Range("A2").Copy()
$ws.Range("CD33").PasteSpecial(-4122)
$ws.Range("CD34").Value = 137.3
$ws.Range("N2").Copy()
$ws.Range("CD34").PasteSpecial(-4122)
$ws.Range("CD35").Value = 127.9
$ws.Range("N2").Copy()
$ws.Range("CD35").PasteSpecial(-4122)
$ws.Range("CD36").Value = 146.1
$ws.Range("A2").Copy()
$ws.Range("CD36").PasteSpecial(-4122)
$ws.Range("CD37").Value = 107.6
$ws.Range("D2").Copy()
$ws.Range("CD37").PasteSpecial(-4122)
$ws.Range("CD38").Value = 181.7
$ws.Range("A2").Copy()
$ws.Range("CD38").PasteSpecial(-4122)
$ws.Range("CD39").Value = 191.8
$ws.Range("A2").Copy()
$ws.Range("CD39").PasteSpecial(-4122)
$ws.Range("CD40").Value = 125.7
$ws.Range("N2").Copy()
$ws.Range("CD40").PasteSpecial(-4122)
$ws.Range("CD41").Value = 150.7
$ws.Range("A2").Copy()
$ws.Range("CD41").PasteSpecial(-4122)
$ws.Range("CD42").Value = 141.7
$ws.Range("A2").Copy()
$ws.Range("CD42").PasteSpecial(-4122)
$ws.Range("CD43").Value = 142.3
$ws.Range("A2").Copy()
$ws.Range("CD43").PasteSpecial(-4122)
$ws.Range("CD44").Value = 152.8
$ws.Range("A2").Copy()
$ws.Range("CD44").PasteSpecial(-4122)
$ws.Range("CD45").Value = 207.7
$ws.Range("A2").Copy()
$ws.Range("CD45").PasteSpecial(-4122)
$ws.Range("CD46").Value = 129.3
$ws.Range("N2").Copy()
$ws.Range("CD46").PasteSpecial(-4122)
$ws.Range("CD47").Value = 195.2
$ws.Range("A2").Copy()
$ws.Range("CD47").PasteSpecial(-4122)
$ws.Range("CD48").Value = 127.6
$ws.Range("N2").Copy()
$ws.Range("CD48").PasteSpecial(-4122)
$ws.Range("CD49").Value = 154.4
$ws.Range("A2").Copy()
$ws.Range("CD49").PasteSpecial(-4122)
$ws.Range("CD50").Value = 153.5
$ws.Range("A2").Copy()
$ws.Range("CD50").PasteSpecial(-4122)
$ws.Range("CD51").Value = 215.5
$ws.Range("A2").Copy()
$ws.Range("CD51").PasteSpecial(-4122)
$ws.Range("CD52").Value = 119.9
$ws.Range("D2").Copy()
$ws.Range("CD52").PasteSpecial(-4122)
$ws.Range("CD53").Value = 141.3
$ws.Range("A2").Copy()
$ws.Range("CD53").PasteSpecial(-4122)
